# "Coupling Parameters" is the 2nd sheet (and the active one) in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) End Year: 2040 -> 2026
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 2026

# ---------------------------------------------------------------------------
# 2) Insert a new row for "Power_plants_from_year" above the old row 5
#    ("CurrentYear"). The threaded comment that lives on B5 ("CurrentYear"
#    value cell) needs to move down to B6 together with the cell content.
# ---------------------------------------------------------------------------
$oldCommentText = $null
$ctList = $ws.CommentsThreaded
if ($ctList.Count -ge 1) {
    $existing = $ctList.Item(1)
    $oldCommentText = $existing.Text()
    $existing.Delete()
}

$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "Power_plants_from_year"
$ws.Range("B5").Value = 2019
$ws.Range("B5").Style = $ws.Range("B4").Style

if ($oldCommentText -ne $null) {
    $ws.Range("B6").AddCommentThreaded($oldCommentText) | Out-Null
}

# ---------------------------------------------------------------------------
# 3) Rename / restyle the (now shifted) fuel-trend and dismantling rows
#    row 12 : start_year_fuel_trends  -> start_tick_fuel_trends
#    row 13 : start_year_dismantling  -> start_tick_dismantling (value + text + style)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "start_tick_fuel_trends"

$ws.Range("A13").Value = "start_tick_dismantling"
$ws.Range("B13").Value = 40
$ws.Range("C13").Value = "Year when the dismantling begins, based on the profits. If this is very high then no dismantling is pconsidered, either for the prepare market clearing"
$ws.Range("C13").Style = $ws.Range("A9").Style
$ws.Rows.Item(13).RowHeight = 29

# ---------------------------------------------------------------------------
# 4) Insert a new row for "install_at_look_ahead_year" above the old blank
#    separator row (currently row 25, right after "writeALLcostsinOPEX").
# ---------------------------------------------------------------------------
$ws.Rows.Item(26).Insert()
$ws.Range("A26").Value = "install_at_look_ahead_year"
$ws.Range("B26").Value = $true
$ws.Rows.Item(26).RowHeight = 13.5
$ws.Rows.Item(27).RowHeight = 13.5

# ---------------------------------------------------------------------------
# 5) Final selection, to mimic the saved cursor position in the workbook.
# ---------------------------------------------------------------------------
$ws.Range("C8").Select()
